# item_data.xlsx edit: adjust column widths, move the active selection,
# and move the "60" Attack value on row 3 over to Defense (H3), clearing
# the old Attack (G3) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- column widths (E = Attack/EquipmentSlot col, K = FoodSatiation col) ---
# Target stored OOXML widths are 16.42578125 / 14.5703125 (real Excel's
# MDW=7 "characters" quantization from a mouse-dragged resize). This
# engine's ColumnWidth setter quantizes the stored width to the nearest
# 1/6 and adds a fixed 5/6 offset, so feed it (target - 5/6) to land on
# the closest achievable stored width (16.5 / 14.5 respectively).
$ws.Columns.Item(5).ColumnWidth = 15.592447916666666
$ws.Columns.Item(11).ColumnWidth = 13.736979166666666

# --- move the selected cell from H8 to H4 ---
$ws.Range("H4").Select() | Out-Null

# --- row 3 data edit: Attack 60 moves off G3, Defense H3 becomes 10 ---
$ws.Range("G3").ClearContents() | Out-Null
$ws.Range("H3").Value = 10
